$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a value to be stored as TEXT (matches the original
# inlineStr cell type) even when it looks numeric, by using the
# classic Excel apostrophe text-prefix so values like "242.45" or
# "0.694" are not silently coerced into numbers.
function Set-TextValue($cellRef, $val) {
    $ws.Range($cellRef).Value = "'" + $val
}

Set-TextValue 'D2' '43.959.17'
Set-TextValue 'E2' '  +0.56%  '
Set-TextValue 'D3' '2.355.98'
Set-TextValue 'E3' '  -0.20%  '
Set-TextValue 'E4' '  +0.13%  '
Set-TextValue 'D5' '0.694'
Set-TextValue 'E5' '  +5.26%  '
Set-TextValue 'D6' '242.45'
Set-TextValue 'D7' '77.00'
Set-TextValue 'E7' '  +4.13%  '
Set-TextValue 'E8' '  +0.04%  '
Set-TextValue 'D9' '0.634'
Set-TextValue 'E9' '  +20.35%  '
Set-TextValue 'E10' '  +4.10%  '
Set-TextValue 'D11' '57.34'
Set-TextValue 'E11' '  +0.73%  '
Set-TextValue 'E12' '  +23.17%  '
Set-TextValue 'D13' '7.53'
Set-TextValue 'E13' '  +13.45%  '
Set-TextValue 'E14' '  +1.76%  '
Set-TextValue 'D15' '2.708.42'
Set-TextValue 'E15' '  +0.01%  '
Set-TextValue 'D16' '16.84'
Set-TextValue 'E16' '  +1.90%  '
Set-TextValue 'E17' '  +5.92%  '
Set-TextValue 'D18' '2.360.50'
Set-TextValue 'E18' '  -0.28%  '
Set-TextValue 'D19' '43.869.58'
Set-TextValue 'E19' '  +0.44%  '
Set-TextValue 'E20' '  +2.22%  '
Set-TextValue 'D21' '6.66'
Set-TextValue 'E21' '  +3.22%  '
Set-TextValue 'D22' '77.77'
Set-TextValue 'E22' '  +2.80%  '
Set-TextValue 'D23' '262.19'
Set-TextValue 'E23' '  +4.46%  '
Set-TextValue 'E24' '  -0.01%  '
Set-TextValue 'E25' '  +1.99%  '
Set-TextValue 'D26' '3.64'
Set-TextValue 'E26' '  -5.14%  '
Set-TextValue 'D27' '10.99'
Set-TextValue 'E27' '  +7.84%  '
Set-TextValue 'E28' '  +17.28%  '
Set-TextValue 'E29' '  +2.56%  '
Set-TextValue 'D30' '23.16'
Set-TextValue 'E30' '  +2.99%  '
Set-TextValue 'D31' '175.30'
Set-TextValue 'E31' '  +1.77%  '
Set-TextValue 'E32' '  -3.31%  '
Set-TextValue 'E33' '  +4.56%  '
Set-TextValue 'D34' '5.39'
Set-TextValue 'D35' '0.0765'
Set-TextValue 'E35' '  +8.98%  '
Set-TextValue 'E36' '  +6.89%  '
Set-TextValue 'E37' '  +2.08%  '
Set-TextValue 'D38' '2.43'
Set-TextValue 'E38' '  -0.22%  '
Set-TextValue 'D39' '6.44'
Set-TextValue 'E39' '  -2.26%  '
Set-TextValue 'D40' '0.0281'
Set-TextValue 'E40' '  +7.53%  '
Set-TextValue 'D41' '0.214'
Set-TextValue 'E41' '  +20.87%  '
Set-TextValue 'D42' '19.40'
Set-TextValue 'E42' '  -0.63%  '
Set-TextValue 'D43' '9.19'
Set-TextValue 'E43' '  +3.62%  '
Set-TextValue 'D44' '0.106'
Set-TextValue 'E44' '  +9.36%  '
Set-TextValue 'E45' '  +0.13%  '
Set-TextValue 'E46' '  +11.32%  '
Set-TextValue 'E47' '  +4.39%  '
$ws.Range('B48').Value = 'Aave'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue 'D48' '102.37'
Set-TextValue 'E48' '  +2.19%  '
$ws.Range('B49').Value = 'ARBITRUM'
$ws.Range('C49').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue 'D49' '1.19'
Set-TextValue 'E49' '  +1.90%  '
Set-TextValue 'D50' '4.56'
Set-TextValue 'E50' '  +2.64%  '
Set-TextValue 'D51' '56.09'
Set-TextValue 'E51' '  +10.11%  '
